$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update footer timestamp (row 1) ---
$ws.Range("A1").Value = 'Datos actualizados a 23 de Marzo de 2020 a las 19:46'

# --- Re-sequence country names (column A) to match updated country list ---
$ws.Range("A28").Value = 'Irlanda'
$ws.Range("A29").Value = 'Japon'
$ws.Range("A30").Value = 'Ecuador'
$ws.Range("A81").Value = 'Jordania'
$ws.Range("A82").Value = 'Vietnam'
$ws.Range("A83").Value = 'Islas Feroe'
$ws.Range("A84").Value = 'Republica de Chipre'
$ws.Range("A120").Value = 'Montenegro'
$ws.Range("A121").Value = 'Costa de Marfil'
$ws.Range("A122").Value = 'Macao'
$ws.Range("A123").Value = 'Ghana'
$ws.Range("A124").Value = 'Monaco'
$ws.Range("A125").Value = 'Paraguay'
$ws.Range("A149").Value = 'Haiti'
$ws.Range("A150").Value = 'Bermudas'
$ws.Range("A154").Value = 'Guinea'
$ws.Range("A155").Value = 'Suazilandia'
$ws.Range("A156").Value = 'Bahamas'
$ws.Range("A159").Value = 'Congo'
$ws.Range("A160").Value = 'San Bartolome'
$ws.Range("A161").Value = 'Namibia'
$ws.Range("A163").Value = 'Zambia'
$ws.Range("A164").Value = 'Cabo Verde'
$ws.Range("A165").Value = 'Fiyi'
$ws.Range("A167").Value = 'Santa Lucia'
$ws.Range("A168").Value = 'Liberia'
$ws.Range("A171").Value = 'Angola'
$ws.Range("A172").Value = 'Birmania'
$ws.Range("A173").Value = 'Nicaragua'
$ws.Range("A174").Value = 'San Martin (Parte Holandesa)'
$ws.Range("A175").Value = 'Butan'
$ws.Range("A176").Value = 'Mauritania'
$ws.Range("A177").Value = 'Benin'
$ws.Range("A178").Value = 'Niger'
$ws.Range("A183").Value = 'Uganda'
$ws.Range("A184").Value = 'Eritrea'
$ws.Range("A185").Value = 'Belice'
$ws.Range("A186").Value = 'Dominica'
$ws.Range("A187").Value = 'San Vicente y las Granadinas'
$ws.Range("A188").Value = 'Timor Oriental'
$ws.Range("A189").Value = 'Papua Nueva Guinea'
$ws.Range("A190").Value = 'Santa Sede'
$ws.Range("A191").Value = 'Somalia'
$ws.Range("A192").Value = 'Republica de Yibuti'
$ws.Range("A193").Value = 'Republica del Chad'
$ws.Range("A194").Value = 'Antigua y Barbuda'
$ws.Range("A195").Value = 'Siria'
$ws.Range("A196").Value = 'Islas Turcas y Caicos'
$ws.Range("A197").Value = 'Mozambique'

# --- Update statistic values (columns B-H) for rows with new data ---
$ws.Range("B6").Value = 41569
$ws.Range("C6").Value = 8023
$ws.Range("E6").Value = 40878
$ws.Range("G6").Value = 85
$ws.Range("H6").Value = 504
$ws.Range("B8").Value = 29056
$ws.Range("C8").Value = 4183
$ws.Range("E8").Value = 28516
$ws.Range("B10").Value = 19856
$ws.Range("C10").Value = 3838
$ws.Range("E10").Value = 16796
$ws.Range("F10").Value = 2082
$ws.Range("B20").Value = 2035
$ws.Range("C20").Value = 565
$ws.Range("D20").Value = 320
$ws.Range("E20").Value = 1692
$ws.Range("B28").Value = 1125
$ws.Range("C28").Value = 219
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 1114
$ws.Range("F28").Value = 29
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 6
$ws.Range("B29").Value = 1101
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 235
$ws.Range("E29").Value = 825
$ws.Range("F29").Value = 49
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 41
$ws.Range("B30").Value = 981
$ws.Range("C30").Value = 192
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = 960
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 18
$ws.Range("E41").Value = 496
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 7
$ws.Range("D76").Value = 5
$ws.Range("E76").Value = 134
$ws.Range("B80").Value = 131
$ws.Range("C80").Value = 5
$ws.Range("E80").Value = 128
$ws.Range("B81").Value = 127
$ws.Range("C81").Value = 15
$ws.Range("D81").Value = 1
$ws.Range("E81").Value = 126
$ws.Range("F81").Value = 0
$ws.Range("B82").Value = 123
$ws.Range("C82").Value = 10
$ws.Range("D82").Value = 17
$ws.Range("E82").Value = 106
$ws.Range("F82").Value = 2
$ws.Range("B83").Value = 118
$ws.Range("C83").Value = 3
$ws.Range("D83").Value = 14
$ws.Range("E83").Value = 104
$ws.Range("F83").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("B84").Value = 116
$ws.Range("C84").Value = 21
$ws.Range("D84").Value = 3
$ws.Range("E84").Value = 112
$ws.Range("F84").Value = 3
$ws.Range("H84").Value = 1
$ws.Range("B120").Value = 27
$ws.Range("C120").Value = 6
$ws.Range("D120").Value = 0
$ws.Range("E120").Value = 26
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 1
$ws.Range("C121").Value = 11
$ws.Range("D121").Value = 2
$ws.Range("E121").Value = 23
$ws.Range("B122").Value = 25
$ws.Range("C122").Value = 3
$ws.Range("D122").Value = 10
$ws.Range("E122").Value = 15
$ws.Range("H122").Value = 0
$ws.Range("B123").Value = 24
$ws.Range("C123").Value = 1
$ws.Range("D123").Value = 0
$ws.Range("E123").Value = 23
$ws.Range("H123").Value = 1
$ws.Range("B124").Value = 23
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 1
$ws.Range("E124").Value = 22
$ws.Range("F124").Value = 0
$ws.Range("H124").Value = 0
$ws.Range("C125").Value = 0
$ws.Range("F125").Value = 1
$ws.Range("G125").Value = 0
$ws.Range("C149").Value = 4
$ws.Range("C150").Value = 0
$ws.Range("C154").Value = 2
$ws.Range("C155").Value = 0
$ws.Range("C163").Value = 0
$ws.Range("C164").Value = 0
$ws.Range("C165").Value = 1
$ws.Range("C167").Value = 1
$ws.Range("C172").Value = 2
$ws.Range("C174").Value = 1
$ws.Range("C175").Value = 0
$ws.Range("C176").Value = 0
$ws.Range("C186").Value = 0
$ws.Range("C196").Value = 1
